$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Estimated")
$ws.Range("B2").Value = "Tue Nov 12 18:09:09 EST 2024"
$ws.Range("B3").Value = "Tue Nov 12 18:09:30 EST 2024"
$ws.Range("B4").Value = "Tue Nov 12 18:09:56 EST 2024"
$ws.Range("B5").Value = "Tue Nov 12 18:10:18 EST 2024"
$ws.Range("B6").Value = "Tue Nov 12 18:10:40 EST 2024"
$ws.Range("B7").Value = "Tue Nov 12 18:11:00 EST 2024"

$ws = $wb.Worksheets.Item("Existing")
$ws.Range("B2").Value = "Tue Nov 12 18:11:19 EST 2024"
$ws.Range("B3").Value = "Tue Nov 12 18:11:40 EST 2024"
$ws.Range("B4").Value = "Tue Nov 12 18:11:58 EST 2024"
$ws.Range("B5").Value = "Tue Nov 12 18:12:16 EST 2024"
$ws.Range("B6").Value = "Tue Nov 12 18:12:38 EST 2024"
$ws.Range("B7").Value = "Tue Nov 12 18:12:56 EST 2024"
$ws.Range("B8").Value = "Tue Nov 12 18:13:16 EST 2024"
$ws.Range("B9").Value = "Tue Nov 12 18:13:37 EST 2024"
$ws.Range("B10").Value = "Tue Nov 12 18:13:59 EST 2024"
$ws.Range("B11").Value = "Tue Nov 12 18:14:18 EST 2024"
$ws.Range("B12").Value = "Tue Nov 12 18:14:37 EST 2024"
$ws.Range("B13").Value = "Tue Nov 12 18:14:58 EST 2024"
$ws.Range("B14").Value = "Tue Nov 12 18:15:18 EST 2024"
$ws.Range("B15").Value = "Tue Nov 12 18:15:39 EST 2024"
$ws.Range("B16").Value = "Tue Nov 12 18:15:58 EST 2024"
$ws.Range("B17").Value = "Tue Nov 12 18:16:19 EST 2024"
$ws.Range("B18").Value = "Tue Nov 12 18:16:41 EST 2024"
$ws.Range("B19").Value = "Tue Nov 12 18:16:59 EST 2024"

$ws = $wb.Worksheets.Item("Extension")
$ws.Range("B2").Value = "Tue Nov 12 18:17:18 EST 2024"
$ws.Range("B3").Value = "Tue Nov 12 18:17:38 EST 2024"
$ws.Range("B4").Value = "Tue Nov 12 18:17:57 EST 2024"
$ws.Range("B5").Value = "Tue Nov 12 18:18:15 EST 2024"
$ws.Range("B6").Value = "Tue Nov 12 18:18:33 EST 2024"
$ws.Range("B7").Value = "Tue Nov 12 18:18:54 EST 2024"

$ws = $wb.Worksheets.Item("NewTaxReturn")
$ws.Range("B2").Value = "Tue Nov 12 18:19:13 EST 2024"
$ws.Range("B3").Value = "Tue Nov 12 18:19:32 EST 2024"
$ws.Range("B4").Value = "Tue Nov 12 18:19:55 EST 2024"
$ws.Range("B5").Value = "Tue Nov 12 18:20:14 EST 2024"
$ws.Range("B6").Value = "Tue Nov 12 18:20:34 EST 2024"
$ws.Range("B7").Value = "Tue Nov 12 18:20:57 EST 2024"
$ws.Range("B8").Value = "Tue Nov 12 18:21:16 EST 2024"
$ws.Range("B9").Value = "Tue Nov 12 18:21:34 EST 2024"
$ws.Range("B10").Value = "Tue Nov 12 18:21:57 EST 2024"
$ws.Range("B11").Value = "Tue Nov 12 18:22:16 EST 2024"
$ws.Range("B12").Value = "Tue Nov 12 18:22:39 EST 2024"
$ws.Range("B13").Value = "Tue Nov 12 18:22:58 EST 2024"
$ws.Range("B14").Value = "Tue Nov 12 18:23:17 EST 2024"
$ws.Range("B15").Value = "Tue Nov 12 18:23:35 EST 2024"
$ws.Range("B16").Value = "Tue Nov 12 18:23:58 EST 2024"
$ws.Range("B17").Value = "Tue Nov 12 18:24:17 EST 2024"
$ws.Range("B18").Value = "Tue Nov 12 18:24:38 EST 2024"
$ws.Range("B19").Value = "Tue Nov 12 18:24:59 EST 2024"
$ws.Range("B20").Value = "Tue Nov 12 18:25:18 EST 2024"
$ws.Range("B21").Value = "Tue Nov 12 18:25:39 EST 2024"
$ws.Range("B22").Value = "Tue Nov 12 18:26:00 EST 2024"
$ws.Range("B23").Value = "Tue Nov 12 18:26:21 EST 2024"
$ws.Range("B24").Value = "Tue Nov 12 18:26:42 EST 2024"
$ws.Range("B25").Value = "Tue Nov 12 18:27:02 EST 2024"
$ws.Range("B26").Value = "Tue Nov 12 18:27:23 EST 2024"
$ws.Range("B27").Value = "Tue Nov 12 18:27:44 EST 2024"
$ws.Range("B28").Value = "Tue Nov 12 18:28:02 EST 2024"
$ws.Range("B29").Value = "Tue Nov 12 18:28:23 EST 2024"
$ws.Range("B30").Value = "Tue Nov 12 18:28:45 EST 2024"
$ws.Range("B31").Value = "Tue Nov 12 18:29:05 EST 2024"
$ws.Range("B32").Value = "Tue Nov 12 18:29:24 EST 2024"
$ws.Range("B33").Value = "Tue Nov 12 18:29:44 EST 2024"
$ws.Range("B34").Value = "Tue Nov 12 18:30:04 EST 2024"
$ws.Range("B35").Value = "Tue Nov 12 18:30:25 EST 2024"
$ws.Range("B36").Value = "Tue Nov 12 18:30:45 EST 2024"
$ws.Range("B37").Value = "Tue Nov 12 18:31:05 EST 2024"
$ws.Range("B38").Value = "Tue Nov 12 18:31:26 EST 2024"
$ws.Range("B39").Value = "Tue Nov 12 18:31:47 EST 2024"
$ws.Range("B40").Value = "Tue Nov 12 18:32:07 EST 2024"
$ws.Range("B41").Value = "Tue Nov 12 18:32:30 EST 2024"
$ws.Range("B42").Value = "Tue Nov 12 18:32:51 EST 2024"
$ws.Range("B43").Value = "Tue Nov 12 18:33:11 EST 2024"
$ws.Range("B44").Value = "Tue Nov 12 18:33:30 EST 2024"
$ws.Range("B45").Value = "Tue Nov 12 18:33:51 EST 2024"
$ws.Range("B46").Value = "Tue Nov 12 18:34:11 EST 2024"
$ws.Range("B47").Value = "Tue Nov 12 18:34:32 EST 2024"
$ws.Range("B48").Value = "Tue Nov 12 18:34:54 EST 2024"
$ws.Range("B49").Value = "Tue Nov 12 18:35:14 EST 2024"
$ws.Range("B50").Value = "Tue Nov 12 18:35:34 EST 2024"
$ws.Range("B51").Value = "Tue Nov 12 18:35:58 EST 2024"
$ws.Range("B52").Value = "Tue Nov 12 18:36:18 EST 2024"

$ws = $wb.Worksheets.Item("Personal_IND")
$ws.Range("B2").Value = "Tue Nov 12 18:37:01 EST 2024"
$ws.Range("B3").Value = "Tue Nov 12 18:37:20 EST 2024"
$ws.Range("B4").Value = "Tue Nov 12 18:37:40 EST 2024"
$ws.Range("B5").Value = "Tue Nov 12 18:37:59 EST 2024"
$ws.Range("B6").Value = "Tue Nov 12 18:38:17 EST 2024"

$ws = $wb.Worksheets.Item("Personal_JNT")
$ws.Range("B2").Value = "Tue Nov 12 18:38:41 EST 2024"
$ws.Range("B3").Value = "Tue Nov 12 18:39:07 EST 2024"
$ws.Range("B4").Value = "Tue Nov 12 18:39:32 EST 2024"
$ws.Range("B5").Value = "Tue Nov 12 18:40:02 EST 2024"
$ws.Range("B6").Value = "Tue Nov 12 18:40:28 EST 2024"

$ws = $wb.Worksheets.Item("Personal_EL")
$ws.Range("B2").Value = "Tue Nov 12 18:36:41 EST 2024"
